$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Add the "scaled" copy of the station coordinate block (C2:F7) into I2:L7 ---
# Copy values first
$vals = $ws.Range("C2:F7").Value2
$ws.Range("I2:L7").Value2 = $vals

# Copy formats cell-by-cell (column E -> K, column F -> L) so existing styles are reused
for ($r = 2; $r -le 7; $r++) {
  $ws.Range("E$r").Copy()
  $ws.Range("K$r").PasteSpecial(-4122)
  $ws.Range("F$r").Copy()
  $ws.Range("L$r").PasteSpecial(-4122)
}

# --- 2. Consolidate duplicate cell styles on E:F (many near-identical xfs collapse to 2) ---
$ws.Range("F2").Copy()
$ws.Range("E8:F22").PasteSpecial(-4122)
$ws.Range("E24:F24").PasteSpecial(-4122)

$ws.Range("E2").Copy()
$ws.Range("E23:F23").PasteSpecial(-4122)

# --- 3. Clear the redundant direct formatting on columns A/B and on D11:D24 ---
$ws.Range("A1:B24").Style = "Normal"
$ws.Range("D11:D24").Style = "Normal"

# Drop the custom column-width/style definition on columns A:B
$ws.Range("A:B").EntireColumn.ClearFormats()

# --- 4. Update the selection to match the authored workbook ---
$ws.Range("K2").Select()
